$wb = $excel.ActiveWorkbook
$wsLogin = $wb.Worksheets.Item("Login")
$wsRRAA  = $wb.Worksheets.Item("RRAA")
$wsDatos = $wb.Worksheets.Item("Datos")

# ---------------------------------------------------------------------------
# 1) Cell VALUE changes, written in an order that reproduces the target
#    shared-strings table ordering (new/unused strings drop out, new ones
#    are appended in first-write order).
# ---------------------------------------------------------------------------

# Replace "AJENICARAGUA S.A." with "EMPRESA PRUEBA" (frees old string).
$wsRRAA.Range("B2").Value = "EMPRESA PRUEBA"

# Temporarily blank the cell holding "FAIL" so that string is dropped
# before we introduce the new Rol Online strings.
$wsRRAA.Range("F2").Value = ""

# Introduce the new strings, in the exact order they appear in the target
# workbook's shared string table.
$wsRRAA.Range("F1").Value = "RolOnline"
$wsDatos.Range("G5").Value = "Rol Online"
$wsDatos.Range("G6").Value = "Decisor"
$wsDatos.Range("G7").Value = "Autorizado"
$wsDatos.Range("G8").Value = "Invitado"

# Now fill in all remaining cells that reuse already-registered strings.
$wsRRAA.Range("F2").Value = "Autorizado"
$wsRRAA.Range("F3").Value = "Decisor"
$wsRRAA.Range("F4").Value = "Autorizado"
$wsRRAA.Range("B3").Value = "EMPRESA PRUEBA"
$wsRRAA.Range("C3").Value = "DNI"
$wsRRAA.Range("E2").Value = "RepresentanteLegal"
$wsRRAA.Range("E3").Value = "Administrativo"
$wsRRAA.Range("D3").Value = 70622837

$wsLogin.Range("C3").Value = "DNI"
$wsLogin.Range("D3").Value = 7240270

# ---------------------------------------------------------------------------
# 2) Styles: reuse existing style entries (copy / alignment tricks) instead
#    of creating duplicate style records.
# ---------------------------------------------------------------------------

# Style index 3 (bold Consolas, centered) - already used on Login!B1:E3.
$fmtBoldCenter = $wsLogin.Range("B1")
$fmtBoldCenter.Copy()
$wsRRAA.Range("B1:F1").PasteSpecial(-4122)
$wsRRAA.Range("B2:E10").PasteSpecial(-4122)

# Style index 4 (horizontal center, default font) - already used by the
# Login sheet's column-level formatting (columns B:E, style 4).
$fmtCenter = $wsLogin.Range("B100")
$fmtCenter.Copy()
$wsRRAA.Range("B11:E15").PasteSpecial(-4122)

# Style index 5 (thin border, no fill) - already used on Datos!E5:E10.
$fmtBordered = $wsDatos.Range("E5")
$fmtBordered.Copy()
$wsDatos.Range("G5:G8").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Hyperlink on Login!B3 (this also switches its style to the built-in
#    Hyperlink style, matching the target cellXfs/font/cellStyle additions).
# ---------------------------------------------------------------------------
$wsLogin.Hyperlinks.Add($wsLogin.Range("B3"), "https://wappe.movistar.com.pe/#/agente", "/agente")

# ---------------------------------------------------------------------------
# 4) Sheet visibility / activation / selection.
# ---------------------------------------------------------------------------
$wsDatos.Visible = $false

$wsLogin.Range("E2").Select()
$wsDatos.Range("F11").Select()

$wsRRAA.Activate()
$wsRRAA.Range("D18").Select()
